$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.801.35'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '2.701.07'
$ws.Range('E3').Value = '  +1.84%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.47'
$ws.Range('E5').Value = '  +1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.65'
$ws.Range('E6').Value = '  +1.03%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.84%  '
$ws.Range('E9').Value = '  +5.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.06'
$ws.Range('E10').Value = '  +4.67%  '
$ws.Range('E11').Value = '  +0.56%  '
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '30.10'
$ws.Range('E13').Value = '  +3.69%  '
$ws.Range('E14').Value = '  +8.45%  '
$ws.Range('D15').Value = '3.185.93'
$ws.Range('E15').Value = '  +1.93%  '
$ws.Range('D16').Value = '65.670.54'
$ws.Range('E16').Value = '  +1.03%  '
$ws.Range('D17').Value = '2.723.30'
$ws.Range('E17').Value = '  +2.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.74'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '360.14'
$ws.Range('E20').Value = '  +1.58%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.54'
$ws.Range('E21').Value = '  +3.52%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('E23').Value = '  +2.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.87'
$ws.Range('E24').Value = '  +3.65%  '
$ws.Range('E25').Value = '  +11.46%  '
$ws.Range('E26').Value = '  -4.38%  '
$ws.Range('E27').Value = '  +2.98%  '
$ws.Range('E28').Value = '  +3.88%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.33'
$ws.Range('E29').Value = '  +1.89%  '
$ws.Range('E30').Value = '  +4.21%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '540.06'
$ws.Range('E32').Value = '  +3.32%  '
$ws.Range('E33').Value = '  +0.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.68'
$ws.Range('E34').Value = '  +5.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.45'
$ws.Range('E35').Value = '  -3.39%  '
$ws.Range('E36').Value = '  +1.09%  '
$ws.Range('E37').Value = '  +2.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.50'
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('E39').Value = '  -0.99%  '
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '42.84'
$ws.Range('E42').Value = '  +1.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '168.23'
$ws.Range('E43').Value = '  +1.71%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.20'
$ws.Range('E44').Value = '  +2.25%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.64'
$ws.Range('E46').Value = '  +2.89%  '
$ws.Range('E47').Value = '  +2.58%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0266'
$ws.Range('E48').Value = '  +4.47%  '
$ws.Range('E49').Value = '  +1.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.09'
$ws.Range('E50').Value = '  +8.20%  '
$ws.Range('E51').Value = '  -0.25%  '
